# Primary and Foreign Key
# Adds a new "Foreign key" block (rows 27-33) below the existing
# "Primary key" notes on Sheet1, and moves the active selection/view
# to reflect where the author ended up working (around column G / cell C34).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A27").Value = "Foreign key"
$ws.Range("B27").Value = "seek"
$ws.Range("B28").Value = "scan"
$ws.Range("B29").Value = "connect 2 tables "
$ws.Range("B30").Value = "cascade effect"
$ws.Range("B31").Value = "null"
$ws.Range("B32").Value = "default"
$ws.Range("B33").Value = "no actions"
$ws.Range("C33").Value = "ERROR"

# Scroll the view over to column G and land the selection on C34,
# matching where the author left off after typing the new notes.
$ws.Application.ActiveWindow.ScrollColumn = 7
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("C34").Select()
